$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force specific Price cells to Text format so that numeric-looking
# strings (e.g. "519.32") are stored as text, matching the source data,
# instead of being auto-converted to numbers by Excel.
$textCells = @("D5", "D6", "D7", "D10", "D11", "D16", "D19", "D20", "D21", "D24", "D25", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D41", "D42", "D43", "D46", "D48", "D50")
$u = $ws.Range($textCells[0])
for ($i = 1; $i -lt $textCells.Length; $i++) {
    $u = $excel.Union($u, $ws.Range($textCells[$i]))
}
foreach ($area in $u.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = "58.208.62"
$ws.Range("E2").Value = "  -1.12%  "

$ws.Range("D3").Value = "2.480.43"
$ws.Range("E3").Value = "  -0.65%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "519.32"
$ws.Range("E5").Value = "  -2.51%  "

$ws.Range("D6").Value = "134.86"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("E8").Value = "  -1.66%  "

$ws.Range("D9").Value = "2.495.35"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").Value = "0.0989"
$ws.Range("E10").Value = "  -2.36%  "

$ws.Range("D11").Value = "0.157"
$ws.Range("E11").Value = "  -0.73%  "

$ws.Range("E12").Value = "  -0.98%  "

$ws.Range("E13").Value = "  -1.90%  "

$ws.Range("D14").Value = "2.921.81"
$ws.Range("E14").Value = "  -0.56%  "

$ws.Range("D15").Value = "58.122.72"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").Value = "22.20"
$ws.Range("E16").Value = "  -2.16%  "

$ws.Range("E17").Value = "  -2.01%  "

$ws.Range("D18").Value = "2.485.01"
$ws.Range("E18").Value = "  -0.21%  "

$ws.Range("D19").Value = "10.68"

$ws.Range("D20").Value = "4.19"
$ws.Range("E20").Value = "  -1.15%  "

$ws.Range("D21").Value = "321.36"
$ws.Range("E21").Value = "  -0.42%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("D24").Value = "64.28"
$ws.Range("E24").Value = "  -1.09%  "

$ws.Range("D25").Value = "0.412"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  -1.08%  "

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.53%  "

$ws.Range("D28").Value = "7.39"
$ws.Range("E28").Value = "  -1.57%  "

$ws.Range("D29").Value = "0.0₃0749"
$ws.Range("E29").Value = "  -1.19%  "

$ws.Range("D30").Value = "169.72"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").Value = "6.33"
$ws.Range("E31").Value = "  -1.60%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  +2.62%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.69"
$ws.Range("E33").Value = "  -2.74%  "

$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("E35").Value = "  -0.19%  "

$ws.Range("D36").Value = "18.12"

$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("D38").Value = "4.03"
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").Value = "36.66"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("E40").Value = "  -2.78%  "

$ws.Range("D41").Value = "0.800"
$ws.Range("E41").Value = "  +0.26%  "

$ws.Range("D42").Value = "5.20"
$ws.Range("E42").Value = "  +4.29%  "

$ws.Range("D43").Value = "276.41"
$ws.Range("E43").Value = "  -1.53%  "

$ws.Range("E44").Value = "  -3.06%  "

$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("D46").Value = "123.96"
$ws.Range("E46").Value = "  -4.20%  "

$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("D48").Value = "0.0492"
$ws.Range("E48").Value = "  -1.23%  "

$ws.Range("E49").Value = "  -1.50%  "

$ws.Range("D50").Value = "17.07"
$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("D51").Value = "1.741.38"
$ws.Range("E51").Value = "  -0.53%  "
